$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item(1).Name = "summ44455485"
$wb.Worksheets.Item(2).Name = "summ44552852"
$wb.Worksheets.Item(3).Name = "summ44650364"
$wb.Worksheets.Item(4).Name = "summ44747331"
$wb.Worksheets.Item(5).Name = "summ44841619"
$wb.Worksheets.Item(6).Name = "summ44936125"
$wb.Worksheets.Item(7).Name = "summ45036395"
$wb.Worksheets.Item(8).Name = "summ45224728"
$wb.Worksheets.Item(9).Name = "summ45333769"
$wb.Worksheets.Item(10).Name = "summ45430285"
$wb.Worksheets.Item(11).Name = "summ45526804"
$wb.Worksheets.Item(12).Name = "summ45641523"
$wb.Worksheets.Item(13).Name = "summ45764194"
$wb.Worksheets.Item(14).Name = "summ45881236"
$wb.Worksheets.Item(15).Name = "summ45993932"
$wb.Worksheets.Item(16).Name = "summ46114640"
$wb.Worksheets.Item(17).Name = "summ46241092"
$wb.Worksheets.Item(18).Name = "summ46360998"
$wb.Worksheets.Item(19).Name = "summ46465322"
$wb.Worksheets.Item(20).Name = "summ46565342"
$wb.Worksheets.Item(21).Name = "summ46662856"
$wb.Worksheets.Item(22).Name = "summ46755865"
$wb.Worksheets.Item(23).Name = "summ46849377"
$wb.Worksheets.Item(24).Name = "summ46942334"
$wb.Worksheets.Item(25).Name = "summ47036845"
$wb.Worksheets.Item(26).Name = "summ47135572"
$wb.Worksheets.Item(27).Name = "summ47231449"
$wb.Worksheets.Item(28).Name = "summ47331294"
$wb.Worksheets.Item(29).Name = "summ47430807"
$wb.Worksheets.Item(30).Name = "summ47529515"
$wb.Worksheets.Item(31).Name = "summ47626280"
$wb.Worksheets.Item(32).Name = "summ47726257"
$wb.Worksheets.Item(33).Name = "summ47822775"
$wb.Worksheets.Item(34).Name = "summ47915807"
$wb.Worksheets.Item(35).Name = "summ48010327"
$wb.Worksheets.Item(36).Name = "summ48107849"
$wb.Worksheets.Item(37).Name = "summ48201363"
$wb.Worksheets.Item(38).Name = "summ48294874"
$wb.Worksheets.Item(39).Name = "summ48393389"
$wb.Worksheets.Item(40).Name = "summ48491068"
$wb.Worksheets.Item(41).Name = "summ48585575"
$wb.Worksheets.Item(42).Name = "summ48685669"
$wb.Worksheets.Item(43).Name = "summ48779441"
$wb.Worksheets.Item(44).Name = "summ48875803"
$wb.Worksheets.Item(45).Name = "summ48969331"
$wb.Worksheets.Item(46).Name = "summ49065845"
$wb.Worksheets.Item(47).Name = "summ49209423"
$wb.Worksheets.Item(48).Name = "summ49307464"
$wb.Worksheets.Item(49).Name = "summ49401973"
$wb.Worksheets.Item(50).Name = "summ49497596"
